$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume cells keep their original Text format so
# values like "20.10" or "1.001" are not coerced into numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.235.87'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.96%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.854.60'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.53%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '313.84'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.84%  '
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.24%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4641'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.23%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3713'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.10%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07289'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.68%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8869'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.24%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '20.10'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.78%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07852'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -0.36%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.820.30'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.393'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.521'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -0.47%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.17%  '
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -0.38%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008934'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.25%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.71'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -0.50%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '27.263.00'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.086'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.52'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.075.93'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.30%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.952'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  +5.65%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.39'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -0.92%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.40'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.07%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.94'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.24%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.052'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08813'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.73%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.139'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +6.30%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7684'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +5.61%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.167'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.512'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +1.73%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.725'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +10.53%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +3.47%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01939'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05223'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.02%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -0.38%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.048'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.84%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5122'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.72%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1628'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +0.29%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.456'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.58%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4795'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -0.80%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.36'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +1.96%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9997'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -0.31%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.99'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.56%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.644'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.96%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06205'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.10%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '65.59'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +1.04%  '
